$wb = $excel.ActiveWorkbook

# Update "展览" (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2553
$ws1.Range("F6").Value = 236

# Update "全部类型" (All Types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2553
$ws4.Range("F6").Value = 236
